$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "RdbmsUUID" header in A1 is renamed to "FullyQualifiedName" (a new
# mapping model entry replacing the old one in the shared-string table).
$ws.Range("A1").Value = "FullyQualifiedName"

# Move/restore the active selection to A2 (was D7).
$ws.Range("A2").Select()
